$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.19
$ws.Range("P2").Value = 4.33
$ws.Range("S2").Value = 1.32
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.87
